$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the sheet's old synthetic "entire column" sentinel row first (while
# it is still at its original address) so the used range tightens back up
# to the real data once the new rows below are added.
$null = $ws.Rows("1048576:1048576").Delete()

# A new row ("location of the base forecast software") is inserted right
# after the header row, pushing the existing rows (client location,
# scenario name, v_date, output location) down by one.
$null = $ws.Rows("2:2").Insert()

$ws.Range("A2").Value = "מיקום תוכנת תחזית בסיס"
$ws.Range("B2").Value = "C:\Users\dpere\Documents\JTMT\forecast_git\create_forecast_basic\current"

# Update the forecast-scenario name (was "with_project").
$ws.Range("B4").Value = "realy"

# Update v_date - now a descriptive string instead of a bare date number.
$ws.Range("B5").Value = "240818_with_poten"

# Preserve the word-wrapped label formatting on these two label cells
# (they keep the same "wrap" style they had before the insert shifted them).
$ws.Range("A4").WrapText = $true
$ws.Range("A5").WrapText = $true

# Row heights: the new row auto-sizes to its own text, and the output-path
# row keeps the shorter height it already had before the insert.
$ws.Rows("2:2").RowHeight = 13.8
$ws.Rows("6:6").RowHeight = 12.75

# Append a new trailing row with the version number.
$ws.Range("A7").Value = "מספר גירסא"
$ws.Range("B7").Value = 1

# Match the author's final selection (cell B6 was the active cell on save).
$null = $ws.Range("B6").Select()
